$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for the 5 members added in "phân công 7" (rows 4-8)
# MSSV (student id) in column C, Họ và tên (full name) in column D,
# Phân công 1..7 (columns E..K) marked as fully completed (1 = 100%).
# The MSSV/name text cells are written in the exact order they were typed
# (matches the shared-string table ordering of the original edit: rows
# 4-6 filled MSSV-then-name, then row 7's MSSV, row 8's MSSV, row 7's
# name, row 8's name).
$ws.Cells.Item(4, 3).Value = "0712187"
$ws.Cells.Item(4, 4).Value = "Lý Hoài"
$ws.Cells.Item(5, 3).Value = "0712188"
$ws.Cells.Item(5, 4).Value = "Phan Lê Huỳnh"
$ws.Cells.Item(6, 3).Value = "0712236"
$ws.Cells.Item(6, 4).Value = "Phan Vũ Lâm"
$ws.Cells.Item(7, 3).Value = "0712365"
$ws.Cells.Item(8, 3).Value = "0712381"
$ws.Cells.Item(7, 4).Value = "Nguyễn Hồ Mẫn Sáng"
$ws.Cells.Item(8, 4).Value = "La Minh Tâm"

for ($row = 4; $row -le 8; $row++) {
    for ($col = 5; $col -le 11; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
}

# Apply the percentage number format (0%) to the newly filled assignment cells
$ws.Range("E4:K8").NumberFormat = "0%"

# Column D width widened to fit the longer names, no longer auto "best fit"
# (21.42578125 character-units in the saved file is the closest reachable
# quantized width to the requested ~21.43; the stored width snaps to
# 1/6-character pixel increments, same as native Excel column resizing)
$ws.Columns.Item(4).ColumnWidth = 20.65

# Selection moved to H10 (matches the final cursor position in the edited file)
$ws.Range("H10").Select()
